# Scheduled runner update: refresh computed market/profit figures (columns H-N)
# for a set of Leve rows across all job sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(115, 8).Value = 509.0909
$ws.Cells.Item(115, 9).Value = 530
$ws.Cells.Item(115, 10).Value = 300
$ws.Cells.Item(115, 11).Value = 1590
$ws.Cells.Item(115, 12).Value = 900
$ws.Cells.Item(115, 13).Value = -23
$ws.Cells.Item(115, 14).Value = -4034
$ws.Cells.Item(129, 8).Value = 2561.3125
$ws.Cells.Item(129, 9).Value = 940.1667
$ws.Cells.Item(129, 10).Value = 3534
$ws.Cells.Item(129, 11).Value = 2820.5001
$ws.Cells.Item(129, 12).Value = 10602
$ws.Cells.Item(129, 13).Value = 2179.4999
$ws.Cells.Item(129, 14).Value = -20602
$ws.Cells.Item(131, 8).Value = 909.55554
$ws.Cells.Item(131, 9).Value = 602.2083
$ws.Cells.Item(131, 11).Value = 1806.6249
$ws.Cells.Item(131, 13).Value = 3233.3751
$ws.Cells.Item(132, 8).Value = 2527510.5
$ws.Cells.Item(132, 9).Value = 3269610.2
$ws.Cells.Item(132, 11).Value = 9808830.600000001
$ws.Cells.Item(132, 13).Value = -9806300.600000001
$ws.Cells.Item(137, 8).Value = 886.5484
$ws.Cells.Item(137, 9).Value = 915.8214
$ws.Cells.Item(137, 10).Value = 613.3333
$ws.Cells.Item(137, 11).Value = 2747.4642
$ws.Cells.Item(137, 12).Value = 1839.9999
$ws.Cells.Item(137, 13).Value = -197.4642000000003
$ws.Cells.Item(137, 14).Value = -6939.9999
$ws.Cells.Item(138, 8).Value = 2930.2778
$ws.Cells.Item(138, 9).Value = 686.3333
$ws.Cells.Item(138, 10).Value = 5174.222
$ws.Cells.Item(138, 11).Value = 2058.9999
$ws.Cells.Item(138, 12).Value = 15522.666
$ws.Cells.Item(138, 13).Value = 3081.0001
$ws.Cells.Item(138, 14).Value = -25802.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5864.4355
$ws.Cells.Item(32, 9).Value = 2496.5576
$ws.Cells.Item(32, 11).Value = 2496.5576
$ws.Cells.Item(32, 13).Value = -2209.5576
$ws.Cells.Item(61, 8).Value = 1941.6
$ws.Cells.Item(61, 9).Value = 1760.3334
$ws.Cells.Item(61, 10).Value = 2666.6667
$ws.Cells.Item(61, 11).Value = 1760.3334
$ws.Cells.Item(61, 12).Value = 2666.6667
$ws.Cells.Item(61, 13).Value = -1548.3334
$ws.Cells.Item(61, 14).Value = -3090.6667
$ws.Cells.Item(88, 8).Value = 2150.625
$ws.Cells.Item(88, 9).Value = 1749.75
$ws.Cells.Item(88, 10).Value = 2551.5
$ws.Cells.Item(88, 11).Value = 1749.75
$ws.Cells.Item(88, 12).Value = 2551.5
$ws.Cells.Item(88, 13).Value = -1343.75
$ws.Cells.Item(88, 14).Value = -3363.5
$ws.Cells.Item(91, 8).Value = 2150.625
$ws.Cells.Item(91, 9).Value = 1749.75
$ws.Cells.Item(91, 10).Value = 2551.5
$ws.Cells.Item(91, 11).Value = 1749.75
$ws.Cells.Item(91, 12).Value = 2551.5
$ws.Cells.Item(91, 13).Value = -345.75
$ws.Cells.Item(91, 14).Value = -5359.5
$ws.Cells.Item(102, 8).Value = 1664.4445
$ws.Cells.Item(102, 9).Value = 1705
$ws.Cells.Item(102, 10).Value = 1583.3334
$ws.Cells.Item(102, 11).Value = 1705
$ws.Cells.Item(102, 12).Value = 1583.3334
$ws.Cells.Item(102, 13).Value = -83
$ws.Cells.Item(102, 14).Value = -4827.3334
$ws.Cells.Item(110, 8).Value = 1039.6
$ws.Cells.Item(110, 9).Value = 715.08
$ws.Cells.Item(110, 10).Value = 1850.9
$ws.Cells.Item(110, 11).Value = 715.08
$ws.Cells.Item(110, 12).Value = 1850.9
$ws.Cells.Item(110, 13).Value = 1329.92
$ws.Cells.Item(110, 14).Value = -5940.9
$ws.Cells.Item(122, 8).Value = 1966.64
$ws.Cells.Item(122, 9).Value = 1957.2
$ws.Cells.Item(122, 10).Value = 1980.8
$ws.Cells.Item(122, 11).Value = 5871.6
$ws.Cells.Item(122, 12).Value = 5942.4
$ws.Cells.Item(122, 13).Value = -3421.6
$ws.Cells.Item(122, 14).Value = -10842.4
$ws.Cells.Item(136, 8).Value = 1941.6
$ws.Cells.Item(136, 9).Value = 1760.3334
$ws.Cells.Item(136, 10).Value = 2666.6667
$ws.Cells.Item(136, 11).Value = 5281.0002
$ws.Cells.Item(136, 12).Value = 8000.000100000001
$ws.Cells.Item(136, 13).Value = -2731.0002
$ws.Cells.Item(136, 14).Value = -13100.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1514.6364
$ws.Cells.Item(86, 9).Value = 1533.6666
$ws.Cells.Item(86, 10).Value = 1473.8572
$ws.Cells.Item(86, 11).Value = 1533.6666
$ws.Cells.Item(86, 12).Value = 1473.8572
$ws.Cells.Item(86, 13).Value = -410.6666
$ws.Cells.Item(86, 14).Value = -3719.8572
$ws.Cells.Item(89, 8).Value = 1514.6364
$ws.Cells.Item(89, 9).Value = 1533.6666
$ws.Cells.Item(89, 10).Value = 1473.8572
$ws.Cells.Item(89, 11).Value = 7668.333000000001
$ws.Cells.Item(89, 12).Value = 7369.286
$ws.Cells.Item(89, 13).Value = -2052.333000000001
$ws.Cells.Item(89, 14).Value = -18601.286
$ws.Cells.Item(105, 8).Value = 1750285.2
$ws.Cells.Item(105, 9).Value = 2842688.5
$ws.Cells.Item(105, 10).Value = 2440
$ws.Cells.Item(105, 11).Value = 2842688.5
$ws.Cells.Item(105, 12).Value = 2440
$ws.Cells.Item(105, 13).Value = -2840941.5
$ws.Cells.Item(105, 14).Value = -5934

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3928.8235
$ws.Cells.Item(16, 9).Value = 2153.6365
$ws.Cells.Item(16, 10).Value = 7183.3335
$ws.Cells.Item(16, 11).Value = 2153.6365
$ws.Cells.Item(16, 12).Value = 7183.3335
$ws.Cells.Item(16, 13).Value = -1866.6365
$ws.Cells.Item(16, 14).Value = -7757.3335
$ws.Cells.Item(58, 8).Value = 1455.8462
$ws.Cells.Item(58, 9).Value = 701.3333
$ws.Cells.Item(58, 10).Value = 1855.2941
$ws.Cells.Item(58, 11).Value = 701.3333
$ws.Cells.Item(58, 12).Value = 1855.2941
$ws.Cells.Item(58, 13).Value = -498.3333
$ws.Cells.Item(58, 14).Value = -2261.2941
$ws.Cells.Item(99, 8).Value = 1418.2941
$ws.Cells.Item(99, 9).Value = 1115.2307
$ws.Cells.Item(99, 10).Value = 2403.25
$ws.Cells.Item(99, 11).Value = 1115.2307
$ws.Cells.Item(99, 12).Value = 2403.25
$ws.Cells.Item(99, 13).Value = 382.7692999999999
$ws.Cells.Item(99, 14).Value = -5399.25
$ws.Cells.Item(113, 8).Value = 3928.8235
$ws.Cells.Item(113, 9).Value = 2153.6365
$ws.Cells.Item(113, 10).Value = 7183.3335
$ws.Cells.Item(113, 11).Value = 2153.6365
$ws.Cells.Item(113, 12).Value = 7183.3335
$ws.Cells.Item(113, 13).Value = 16.36349999999993
$ws.Cells.Item(113, 14).Value = -11523.3335
$ws.Cells.Item(122, 8).Value = 1017.25806
$ws.Cells.Item(122, 9).Value = 941.7222
$ws.Cells.Item(122, 10).Value = 1121.8462
$ws.Cells.Item(122, 11).Value = 2825.1666
$ws.Cells.Item(122, 12).Value = 3365.5386
$ws.Cells.Item(122, 13).Value = -375.1666
$ws.Cells.Item(122, 14).Value = -8265.5386
$ws.Cells.Item(126, 8).Value = 1418.2941
$ws.Cells.Item(126, 9).Value = 1115.2307
$ws.Cells.Item(126, 10).Value = 2403.25
$ws.Cells.Item(126, 11).Value = 3345.6921
$ws.Cells.Item(126, 12).Value = 7209.75
$ws.Cells.Item(126, 13).Value = -875.6921000000002
$ws.Cells.Item(126, 14).Value = -12149.75
$ws.Cells.Item(132, 8).Value = 1515.92
$ws.Cells.Item(132, 9).Value = 740.5714
$ws.Cells.Item(132, 10).Value = 2502.7273
$ws.Cells.Item(132, 11).Value = 2221.7142
$ws.Cells.Item(132, 12).Value = 7508.1819
$ws.Cells.Item(132, 13).Value = 308.2857999999997
$ws.Cells.Item(132, 14).Value = -12568.1819
$ws.Cells.Item(134, 8).Value = 2619.7144
$ws.Cells.Item(134, 9).Value = 2106.889
$ws.Cells.Item(134, 10).Value = 3542.8
$ws.Cells.Item(134, 11).Value = 6320.667
$ws.Cells.Item(134, 12).Value = 10628.4
$ws.Cells.Item(134, 13).Value = -3785.667
$ws.Cells.Item(134, 14).Value = -15698.4
$ws.Cells.Item(136, 8).Value = 1455.8462
$ws.Cells.Item(136, 9).Value = 701.3333
$ws.Cells.Item(136, 10).Value = 1855.2941
$ws.Cells.Item(136, 11).Value = 2103.9999
$ws.Cells.Item(136, 12).Value = 5565.8823
$ws.Cells.Item(136, 13).Value = 446.0001000000002
$ws.Cells.Item(136, 14).Value = -10665.8823

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 34.333332
$ws.Cells.Item(2, 9).Value = 36.545456
$ws.Cells.Item(2, 10).Value = 10
$ws.Cells.Item(2, 11).Value = 36.545456
$ws.Cells.Item(2, 12).Value = 10
$ws.Cells.Item(2, 13).Value = 76.454544
$ws.Cells.Item(2, 14).Value = -236
$ws.Cells.Item(70, 8).Value = 4903.4194
$ws.Cells.Item(70, 9).Value = 4679.3687
$ws.Cells.Item(70, 10).Value = 5258.1665
$ws.Cells.Item(70, 11).Value = 4679.3687
$ws.Cells.Item(70, 12).Value = 5258.1665
$ws.Cells.Item(70, 13).Value = -4409.3687
$ws.Cells.Item(70, 14).Value = -5798.1665
$ws.Cells.Item(73, 8).Value = 4903.4194
$ws.Cells.Item(73, 9).Value = 4679.3687
$ws.Cells.Item(73, 10).Value = 5258.1665
$ws.Cells.Item(73, 11).Value = 4679.3687
$ws.Cells.Item(73, 12).Value = 5258.1665
$ws.Cells.Item(73, 13).Value = -3743.3687
$ws.Cells.Item(73, 14).Value = -7130.1665
$ws.Cells.Item(102, 8).Value = 2892.5
$ws.Cells.Item(102, 9).Value = 1936.7778
$ws.Cells.Item(102, 10).Value = 4612.8
$ws.Cells.Item(102, 11).Value = 1936.7778
$ws.Cells.Item(102, 12).Value = 4612.8
$ws.Cells.Item(102, 13).Value = -314.7778000000001
$ws.Cells.Item(102, 14).Value = -7856.8
$ws.Cells.Item(113, 8).Value = 10178.833
$ws.Cells.Item(113, 9).Value = 14040
$ws.Cells.Item(113, 11).Value = 14040
$ws.Cells.Item(113, 13).Value = -11870
$ws.Cells.Item(122, 8).Value = 1687.9445
$ws.Cells.Item(122, 9).Value = 1073.375
$ws.Cells.Item(122, 10).Value = 2179.6
$ws.Cells.Item(122, 11).Value = 3220.125
$ws.Cells.Item(122, 12).Value = 6538.799999999999
$ws.Cells.Item(122, 13).Value = -770.125
$ws.Cells.Item(122, 14).Value = -11438.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2105.5557
$ws.Cells.Item(7, 9).Value = 2118.75
$ws.Cells.Item(7, 10).Value = 2000
$ws.Cells.Item(7, 11).Value = 2118.75
$ws.Cells.Item(7, 12).Value = 2000
$ws.Cells.Item(7, 13).Value = -2006.75
$ws.Cells.Item(7, 14).Value = -2224
$ws.Cells.Item(40, 8).Value = 2273.389
$ws.Cells.Item(40, 9).Value = 2276
$ws.Cells.Item(40, 10).Value = 2252.5
$ws.Cells.Item(40, 11).Value = 2276
$ws.Cells.Item(40, 12).Value = 2252.5
$ws.Cells.Item(40, 13).Value = -2140
$ws.Cells.Item(40, 14).Value = -2524.5
$ws.Cells.Item(126, 8).Value = 2105.5557
$ws.Cells.Item(126, 9).Value = 2118.75
$ws.Cells.Item(126, 10).Value = 2000
$ws.Cells.Item(126, 11).Value = 6356.25
$ws.Cells.Item(126, 12).Value = 6000
$ws.Cells.Item(126, 13).Value = -3886.25
$ws.Cells.Item(126, 14).Value = -10940

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3483.3333
$ws.Cells.Item(62, 9).Value = 2250
$ws.Cells.Item(62, 10).Value = 4100
$ws.Cells.Item(62, 11).Value = 2250
$ws.Cells.Item(62, 12).Value = 4100
$ws.Cells.Item(62, 13).Value = -1626
$ws.Cells.Item(62, 14).Value = -5348
$ws.Cells.Item(65, 8).Value = 3483.3333
$ws.Cells.Item(65, 9).Value = 2250
$ws.Cells.Item(65, 10).Value = 4100
$ws.Cells.Item(65, 11).Value = 11250
$ws.Cells.Item(65, 12).Value = 20500
$ws.Cells.Item(65, 13).Value = -8130
$ws.Cells.Item(65, 14).Value = -26740
$ws.Cells.Item(136, 8).Value = 5271.625
$ws.Cells.Item(136, 9).Value = 1107.8
$ws.Cells.Item(136, 10).Value = 12211.333
$ws.Cells.Item(136, 11).Value = 3323.4
$ws.Cells.Item(136, 12).Value = 36633.999
$ws.Cells.Item(136, 13).Value = -773.3999999999996
$ws.Cells.Item(136, 14).Value = -41733.999
